# Lancers ("案件情報.xlsx" / sheet "ランサーズ") scrape refresh:
#   - timestamp bumped to 2025-10-11 06:24:35
#   - data rows trimmed from 19 (rows 2-20) down to 6 (rows 2-7)
#   - new scrape content written into rows 2-7
#   - hyperlinks rebuilt to match the new URLs
#   - column B/D widths tweaked

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# 1) Drop the old rows 8-20 entirely (also shrinks the used-range/dimension to H7).
$ws.Rows("8:20").Delete()

# 2) The engine's row-delete does not prune dangling hyperlink relationships
#    pointing at the removed rows, and Range-scoped Hyperlinks.Delete() ends up
#    clearing the whole sheet anyway - so just clear every hyperlink now and
#    re-add the ones we still need (F2:F7) further down.
$ws.Hyperlinks.Delete()

# 3) Rewrite the surviving data rows (2-7) with the refreshed scrape content.
$ws.Range("A2").Value = "2025-10-11 06:24:35"
$ws.Range("B2").Value = "急募 PR Zoom/Meet×TLDV×ChatGPT×Notion×Slack 議事録ワークフロー構築依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5410688"
$ws.Range("G2").Value = 323
$ws.Range("H2").Value = "🔥GPT,ChatGPT"

$ws.Range("A3").Value = "2025-10-11 06:24:35"
$ws.Range("B3").Value = "【急募】配送状況を自動取得するAPI開発者募集!"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5411268"
$ws.Range("G3").Value = 238
$ws.Range("H3").Value = "🔥API ◆開発"

$ws.Range("A4").Value = "2025-10-11 06:24:35"
$ws.Range("B4").Value = "急募バックエンドエンジニア マッチングサイトの開発"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5404059"
$ws.Range("G4").Value = 93
$ws.Range("H4").Value = "◆開発 ◇サイト"

$ws.Range("A5").Value = "2025-10-11 06:24:35"
$ws.Range("B5").Value = "スプレッドシートをもとにした顧客・売上管理アプリのグライド化(Glide/無料版)"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5411304"
$ws.Range("G5").Value = 55
$ws.Range("H5").Value = "◇アプリ"

$ws.Range("A6").Value = "2025-10-11 06:24:35"
$ws.Range("B6").Value = "【急募】時間単位で入札できるシステム構築の依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5411365"
$ws.Range("G6").Value = 33
$ws.Range("H6").ClearContents()

$ws.Range("A7").Value = "2025-10-11 06:24:35"
$ws.Range("B7").Value = "【フォーム制作】物件見積り査定フォーム制作の依頼"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5411435"
$ws.Range("G7").Value = 13
$ws.Range("H7").ClearContents()

# 4) Re-create the hyperlinks for the URL column (F2:F7), in order, so the
#    relationship ids line up with the refreshed rows.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5410688")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5411268")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5404059")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5411304")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5411365")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5411435")

# 5) Column width tweaks (B: 54 -> 57, D: 30 -> 28). The engine stores XML
#    width as ColumnWidth + 5/6, so back that constant out to land on an
#    exact integer in the saved file.
$ws.Columns("B").ColumnWidth = 57 - 5/6
$ws.Columns("D").ColumnWidth = 28 - 5/6
